# Update NATMI LR-pair TPM-derived metrics for Pdgfc-Pdgfrb (new TPM values)
$ws = $excel.ActiveWorkbook.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 0.362096
$ws.Cells.Item(2, 8).Value = 1.086288
$ws.Cells.Item(2, 9).Value = 0.01048575317613816
$ws.Cells.Item(2, 10).Value = 0.01048575317613816
$ws.Cells.Item(2, 13).Value = 8.226724333333333
$ws.Cells.Item(2, 14).Value = 24.680173
$ws.Cells.Item(2, 15).Value = 0.06198126651953669
$ws.Cells.Item(2, 16).Value = 0.06198126651953669
$ws.Cells.Item(2, 17).Value = 2.978863974202667
$ws.Cells.Item(2, 18).Value = 26.809775767824
$ws.Cells.Item(2, 19).Value = 0.0006499202622682976
$ws.Cells.Item(2, 20).Value = 0.0006499202622682976

# Row 3
$ws.Cells.Item(3, 7).Value = 0.362096
$ws.Cells.Item(3, 8).Value = 1.086288
$ws.Cells.Item(3, 9).Value = 0.01048575317613816
$ws.Cells.Item(3, 10).Value = 0.01048575317613816
$ws.Cells.Item(3, 15).Value = 0.6623065855236785
$ws.Cells.Item(3, 16).Value = 0.6623065855236785
$ws.Cells.Item(3, 17).Value = 31.83092792838933
$ws.Cells.Item(3, 18).Value = 286.478351355504
$ws.Cells.Item(3, 19).Value = 0.00694478338273213
$ws.Cells.Item(3, 20).Value = 0.00694478338273213

# Row 4
$ws.Cells.Item(4, 7).Value = 0.362096
$ws.Cells.Item(4, 8).Value = 1.086288
$ws.Cells.Item(4, 9).Value = 0.01048575317613816
$ws.Cells.Item(4, 10).Value = 0.01048575317613816
$ws.Cells.Item(4, 13).Value = 36.43008433333333
$ws.Cells.Item(4, 14).Value = 109.290253
$ws.Cells.Item(4, 15).Value = 0.2744692388979848
$ws.Cells.Item(4, 16).Value = 0.2744692388979848
$ws.Cells.Item(4, 17).Value = 13.19118781676267
$ws.Cells.Item(4, 18).Value = 118.720690350864
$ws.Cells.Item(4, 19).Value = 0.002878016693526767
$ws.Cells.Item(4, 20).Value = 0.002878016693526767

# Row 5
$ws.Cells.Item(5, 7).Value = 0.362096
$ws.Cells.Item(5, 8).Value = 1.086288
$ws.Cells.Item(5, 9).Value = 0.01048575317613816
$ws.Cells.Item(5, 10).Value = 0.01048575317613816
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.1649703333333333
$ws.Cells.Item(5, 14).Value = 0.494911
$ws.Cells.Item(5, 15).Value = 0.00124290905879997
$ws.Cells.Item(5, 16).Value = 0.00124290905879997
$ws.Cells.Item(5, 17).Value = 0.05973509781866666
$ws.Cells.Item(5, 18).Value = 0.5376158803679999
$ws.Cells.Item(5, 19).Value = 0.00001303283761096267
$ws.Cells.Item(5, 20).Value = 0.00001303283761096267

# Row 6
$ws.Cells.Item(6, 9).Value = 0.09151676111574511
$ws.Cells.Item(6, 10).Value = 0.09151676111574511
$ws.Cells.Item(6, 13).Value = 8.226724333333333
$ws.Cells.Item(6, 14).Value = 24.680173
$ws.Cells.Item(6, 15).Value = 0.06198126651953669
$ws.Cells.Item(6, 16).Value = 0.06198126651953669
$ws.Cells.Item(6, 17).Value = 25.99870301580066
$ws.Cells.Item(6, 18).Value = 233.988327142206
$ws.Cells.Item(6, 19).Value = 0.00567232476171977
$ws.Cells.Item(6, 20).Value = 0.00567232476171977

# Row 7
$ws.Cells.Item(7, 9).Value = 0.09151676111574511
$ws.Cells.Item(7, 10).Value = 0.09151676111574511
$ws.Cells.Item(7, 15).Value = 0.6623065855236785
$ws.Cells.Item(7, 16).Value = 0.6623065855236785
$ws.Cells.Item(7, 19).Value = 0.06061215357275529
$ws.Cells.Item(7, 20).Value = 0.06061215357275529

# Row 8
$ws.Cells.Item(8, 9).Value = 0.09151676111574511
$ws.Cells.Item(8, 10).Value = 0.09151676111574511
$ws.Cells.Item(8, 13).Value = 36.43008433333333
$ws.Cells.Item(8, 14).Value = 109.290253
$ws.Cells.Item(8, 15).Value = 0.2744692388979848
$ws.Cells.Item(8, 16).Value = 0.2744692388979848
$ws.Cells.Item(8, 17).Value = 115.1290483364407
$ws.Cells.Item(8, 18).Value = 1036.161435027966
$ws.Cells.Item(8, 19).Value = 0.02511853576984725
$ws.Cells.Item(8, 20).Value = 0.02511853576984725

# Row 9
$ws.Cells.Item(9, 9).Value = 0.09151676111574511
$ws.Cells.Item(9, 10).Value = 0.09151676111574511
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.1649703333333333
$ws.Cells.Item(9, 14).Value = 0.494911
$ws.Cells.Item(9, 15).Value = 0.00124290905879997
$ws.Cells.Item(9, 16).Value = 0.00124290905879997
$ws.Cells.Item(9, 17).Value = 0.5213514552046666
$ws.Cells.Item(9, 18).Value = 4.692163096842
$ws.Cells.Item(9, 19).Value = 0.0001137470114227924
$ws.Cells.Item(9, 20).Value = 0.0001137470114227924

# Row 10
$ws.Cells.Item(10, 7).Value = 21.49951033333333
$ws.Cells.Item(10, 8).Value = 64.498531
$ws.Cells.Item(10, 9).Value = 0.6225933419953967
$ws.Cells.Item(10, 10).Value = 0.6225933419953966
$ws.Cells.Item(10, 13).Value = 8.226724333333333
$ws.Cells.Item(10, 14).Value = 24.680173
$ws.Cells.Item(10, 15).Value = 0.06198126651953669
$ws.Cells.Item(10, 16).Value = 0.06198126651953669
$ws.Cells.Item(10, 17).Value = 176.8705448139848
$ws.Cells.Item(10, 18).Value = 1591.834903325863
$ws.Cells.Item(10, 19).Value = 0.03858912386350574
$ws.Cells.Item(10, 20).Value = 0.03858912386350574

# Row 11
$ws.Cells.Item(11, 7).Value = 21.49951033333333
$ws.Cells.Item(11, 8).Value = 64.498531
$ws.Cells.Item(11, 9).Value = 0.6225933419953967
$ws.Cells.Item(11, 10).Value = 0.6225933419953966
$ws.Cells.Item(11, 15).Value = 0.6623065855236785
$ws.Cells.Item(11, 16).Value = 0.6623065855236785
$ws.Cells.Item(11, 17).Value = 1889.966649496253
$ws.Cells.Item(11, 18).Value = 17009.69984546627
$ws.Cells.Item(11, 19).Value = 0.412347670506747
$ws.Cells.Item(11, 20).Value = 0.412347670506747

# Row 12
$ws.Cells.Item(12, 7).Value = 21.49951033333333
$ws.Cells.Item(12, 8).Value = 64.498531
$ws.Cells.Item(12, 9).Value = 0.6225933419953967
$ws.Cells.Item(12, 10).Value = 0.6225933419953966
$ws.Cells.Item(12, 13).Value = 36.43008433333333
$ws.Cells.Item(12, 14).Value = 109.290253
$ws.Cells.Item(12, 15).Value = 0.2744692388979848
$ws.Cells.Item(12, 16).Value = 0.2744692388979848
$ws.Cells.Item(12, 17).Value = 783.2289745687048
$ws.Cells.Item(12, 18).Value = 7049.060771118344
$ws.Cells.Item(12, 19).Value = 0.1708827207204293
$ws.Cells.Item(12, 20).Value = 0.1708827207204293

# Row 13
$ws.Cells.Item(13, 7).Value = 21.49951033333333
$ws.Cells.Item(13, 8).Value = 64.498531
$ws.Cells.Item(13, 9).Value = 0.6225933419953967
$ws.Cells.Item(13, 10).Value = 0.6225933419953966
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 0.1649703333333333
$ws.Cells.Item(13, 14).Value = 0.494911
$ws.Cells.Item(13, 15).Value = 0.00124290905879997
$ws.Cells.Item(13, 16).Value = 0.00124290905879997
$ws.Cells.Item(13, 17).Value = 3.546781386193445
$ws.Cells.Item(13, 18).Value = 31.921032475741
$ws.Cells.Item(13, 19).Value = 0.0007738269047146262
$ws.Cells.Item(13, 20).Value = 0.000773826904714626

# Row 14
$ws.Cells.Item(14, 7).Value = 9.510307666666666
$ws.Cells.Item(14, 8).Value = 28.530923
$ws.Cells.Item(14, 9).Value = 0.27540414371272
$ws.Cells.Item(14, 10).Value = 0.27540414371272
$ws.Cells.Item(14, 13).Value = 8.226724333333333
$ws.Cells.Item(14, 14).Value = 24.680173
$ws.Cells.Item(14, 15).Value = 0.06198126651953669
$ws.Cells.Item(14, 16).Value = 0.06198126651953669
$ws.Cells.Item(14, 17).Value = 78.23867949885322
$ws.Cells.Item(14, 18).Value = 704.148115489679
$ws.Cells.Item(14, 19).Value = 0.01706989763204289
$ws.Cells.Item(14, 20).Value = 0.01706989763204289

# Row 15
$ws.Cells.Item(15, 7).Value = 9.510307666666666
$ws.Cells.Item(15, 8).Value = 28.530923
$ws.Cells.Item(15, 9).Value = 0.27540414371272
$ws.Cells.Item(15, 10).Value = 0.27540414371272
$ws.Cells.Item(15, 15).Value = 0.6623065855236785
$ws.Cells.Item(15, 16).Value = 0.6623065855236785
$ws.Cells.Item(15, 17).Value = 836.0266832952453
$ws.Cells.Item(15, 18).Value = 7524.240149657208
$ws.Cells.Item(15, 19).Value = 0.1824019780614441
$ws.Cells.Item(15, 20).Value = 0.1824019780614441

# Row 16
$ws.Cells.Item(16, 7).Value = 9.510307666666666
$ws.Cells.Item(16, 8).Value = 28.530923
$ws.Cells.Item(16, 9).Value = 0.27540414371272
$ws.Cells.Item(16, 10).Value = 0.27540414371272
$ws.Cells.Item(16, 13).Value = 36.43008433333333
$ws.Cells.Item(16, 14).Value = 109.290253
$ws.Cells.Item(16, 15).Value = 0.2744692388979848
$ws.Cells.Item(16, 16).Value = 0.2744692388979848
$ws.Cells.Item(16, 17).Value = 346.4613103326132
$ws.Cells.Item(16, 18).Value = 3118.151792993519
$ws.Cells.Item(16, 19).Value = 0.0755899657141815
$ws.Cells.Item(16, 20).Value = 0.0755899657141815

# Row 17
$ws.Cells.Item(17, 7).Value = 9.510307666666666
$ws.Cells.Item(17, 8).Value = 28.530923
$ws.Cells.Item(17, 9).Value = 0.27540414371272
$ws.Cells.Item(17, 10).Value = 0.27540414371272
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 0.1649703333333333
$ws.Cells.Item(17, 14).Value = 0.494911
$ws.Cells.Item(17, 15).Value = 0.00124290905879997
$ws.Cells.Item(17, 16).Value = 0.00124290905879997
$ws.Cells.Item(17, 17).Value = 1.568918625872555
$ws.Cells.Item(17, 18).Value = 14.120267632853
$ws.Cells.Item(17, 19).Value = 0.0003423023050515884
$ws.Cells.Item(17, 20).Value = 0.0003423023050515884

